$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# ---------------------------------------------------------------------------
# 1. Insert a new column "load_lineshape" before the old "exclude" column
#    (old column U "exclude" shifts right to become column V; old column T
#    "trackbg" stays at T and is renamed "track_bg").
# ---------------------------------------------------------------------------
$ws.Columns("U").Insert()

# ---------------------------------------------------------------------------
# 2. Fix bg_freq content for row 4 (2024-08-08_J_e)
# ---------------------------------------------------------------------------
$ws.Range("M4").Value = "[[-4.1, -4.04],[-3.96, -3.9]]"

# ---------------------------------------------------------------------------
# 3. Headers
# ---------------------------------------------------------------------------
$ws.Range("U1").Value = "load_lineshape"
$ws.Range("T1").Value = "track_bg"

# ---------------------------------------------------------------------------
# 4. track_bg (T) used to be blank for rows 2-15 and 1 for rows 16-24;
#    make every data row explicit, and strip the stray style some cells had.
#    load_lineshape (U) mirrors track_bg for every row.
# ---------------------------------------------------------------------------
$ws.Range("T2:U24").ClearFormats()

for ($r = 2; $r -le 15; $r++) {
    $ws.Range("T$r").Value = 0
    $ws.Range("U$r").Value = 0
}
for ($r = 16; $r -le 24; $r++) {
    $ws.Range("T$r").Value = 1
    $ws.Range("U$r").Value = 1
}

# ---------------------------------------------------------------------------
# 5. Append two "smallbox" rows, copied (values+formats) from the
#    2024-10-02_C_e / 2024-10-03_C_e rows, then set the new
#    load_lineshape/exclude columns.
# ---------------------------------------------------------------------------
$ws.Range("A19:V19").Copy($ws.Range("A25:V25"))
$ws.Range("A25").Value = "2024-10-02_C_e_smallbox"
$ws.Range("T25").Value = 1
$ws.Range("U25").Value = 1
$ws.Range("V25").Value = 0

$ws.Range("A20:V20").Copy($ws.Range("A26:V26"))
$ws.Range("A26").Value = "2024-10-03_C_e_smallbox"
$ws.Range("T26").Value = 1
$ws.Range("U26").Value = 1
$ws.Range("V26").Value = 0

$excel.CutCopyMode = 0

# ---------------------------------------------------------------------------
# 6. Selection moved to F13 (and no longer frozen/scrolled to topLeftCell H1)
# ---------------------------------------------------------------------------
$ws.Range("F13").Select()
